$wb = $excel.ActiveWorkbook

# --- Sheet: lsh_unit_categories (row 28, column C value change) ---
$wsUnit = $wb.Worksheets.Item("lsh_unit_categories")
$wsUnit.Cells.Item(28, 3).Value2 = 'inpatient_ward_transfer_geriatric'

# --- Sheet: lsh_unit_categories (new rows 61-81, isolation-derived unit categories) ---
$wsUnit.Cells.Item(61, 1).Value2 = 'Hjartadeild (Hb-14EG)'
$wsUnit.Cells.Item(61, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(61, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(61, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(61, 5).Value2 = 2

$wsUnit.Cells.Item(62, 1).Value2 = 'Heila-, tauga- og bæklunarskurðdeild (Fv-B6)'
$wsUnit.Cells.Item(62, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(62, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(62, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(62, 5).Value2 = 2

$wsUnit.Cells.Item(63, 1).Value2 = 'Móttökugeðdeild (Hb-33A)'
$wsUnit.Cells.Item(63, 2).Value2 = 'Göngudeild'
$wsUnit.Cells.Item(63, 3).Value2 = 'outpatient_clinic'
$wsUnit.Cells.Item(63, 4).Value2 = 'home'
$wsUnit.Cells.Item(63, 5).Value2 = 1

$wsUnit.Cells.Item(64, 1).Value2 = 'Sérhæfð endurhæfingargeðdeild (Kl-10C)'
$wsUnit.Cells.Item(64, 2).Value2 = 'Göngudeild'
$wsUnit.Cells.Item(64, 3).Value2 = 'outpatient_clinic'
$wsUnit.Cells.Item(64, 4).Value2 = 'home'
$wsUnit.Cells.Item(64, 5).Value2 = 1

$wsUnit.Cells.Item(65, 1).Value2 = 'Móttaka bráða- og göngudeildar (Fv-G2/G3)'
$wsUnit.Cells.Item(65, 2).Value2 = 'Bráðamóttaka'
$wsUnit.Cells.Item(65, 3).Value2 = 'emergency_room'
$wsUnit.Cells.Item(65, 4).Value2 = 'home'
$wsUnit.Cells.Item(65, 5).Value2 = 1

$wsUnit.Cells.Item(66, 1).Value2 = 'Móttaka lyf- og skurðlækninga Fossvogi'
$wsUnit.Cells.Item(66, 2).Value2 = 'Göngudeild'
$wsUnit.Cells.Item(66, 3).Value2 = 'outpatient_clinic'
$wsUnit.Cells.Item(66, 4).Value2 = 'home'
$wsUnit.Cells.Item(66, 5).Value2 = 1

$wsUnit.Cells.Item(67, 1).Value2 = 'Bráðageðdeild 32C (Hb-32C)'
$wsUnit.Cells.Item(67, 2).Value2 = 'Bráðamóttaka'
$wsUnit.Cells.Item(67, 3).Value2 = 'emergency_room'
$wsUnit.Cells.Item(67, 4).Value2 = 'home'
$wsUnit.Cells.Item(67, 5).Value2 = 1

$wsUnit.Cells.Item(68, 1).Value2 = 'Útskriftardeild aldraðra (Lk-L2)'
$wsUnit.Cells.Item(68, 2).Value2 = 'Göngudeild'
$wsUnit.Cells.Item(68, 3).Value2 = 'outpatient_clinic'
$wsUnit.Cells.Item(68, 4).Value2 = 'home'
$wsUnit.Cells.Item(68, 5).Value2 = 1

$wsUnit.Cells.Item(69, 1).Value2 = 'Gigtar-og almenn lyflækningadeild (Fv-B7)'
$wsUnit.Cells.Item(69, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(69, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(69, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(69, 5).Value2 = 2

$wsUnit.Cells.Item(70, 1).Value2 = 'HNE-, lýta- og æðaskurðdeild (Fv-A4)'
$wsUnit.Cells.Item(70, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(70, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(70, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(70, 5).Value2 = 2

$wsUnit.Cells.Item(71, 1).Value2 = 'Kvenlækningadeild (Hb-21A)'
$wsUnit.Cells.Item(71, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(71, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(71, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(71, 5).Value2 = 2

$wsUnit.Cells.Item(72, 1).Value2 = 'Bráðalyflækningadeild (Fv-A2)'
$wsUnit.Cells.Item(72, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(72, 3).Value2 = 'inpatient_ward_transfer'
$wsUnit.Cells.Item(72, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(72, 5).Value2 = 2

$wsUnit.Cells.Item(73, 1).Value2 = 'Göngudeild og bráðamóttaka BH (Hb-20E/D)'
$wsUnit.Cells.Item(73, 2).Value2 = 'Bráðamóttaka'
$wsUnit.Cells.Item(73, 3).Value2 = 'emergency_room'
$wsUnit.Cells.Item(73, 4).Value2 = 'home'
$wsUnit.Cells.Item(73, 5).Value2 = 1

$wsUnit.Cells.Item(74, 1).Value2 = 'Taugalækningadeild (Fv-B2)'
$wsUnit.Cells.Item(74, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(74, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(74, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(74, 5).Value2 = 2

$wsUnit.Cells.Item(75, 1).Value2 = 'Endurhæfingardeild (Gr-R2)'
$wsUnit.Cells.Item(75, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(75, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(75, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(75, 5).Value2 = 2

$wsUnit.Cells.Item(76, 1).Value2 = 'Líknardeild, legudeild (Kv-h10-1h)'
$wsUnit.Cells.Item(76, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(76, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(76, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(76, 5).Value2 = 2

$wsUnit.Cells.Item(77, 1).Value2 = 'Dag- og göngudeild blóð- og krabbameinslækninga (H'
$wsUnit.Cells.Item(77, 2).Value2 = 'Dagdeild'
$wsUnit.Cells.Item(77, 3).Value2 = 'outpatient_clinic'
$wsUnit.Cells.Item(77, 4).Value2 = 'home'
$wsUnit.Cells.Item(77, 5).Value2 = 1

$wsUnit.Cells.Item(78, 1).Value2 = 'Blóð- og krabbameinslækningadeild (Hb-11EG)'
$wsUnit.Cells.Item(78, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(78, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(78, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(78, 5).Value2 = 2

$wsUnit.Cells.Item(79, 1).Value2 = 'Dagdeild skurðlækninga F, legudeild (Fv-A5)'
$wsUnit.Cells.Item(79, 2).Value2 = 'Dagdeild'
$wsUnit.Cells.Item(79, 3).Value2 = 'outpatient_clinic'
$wsUnit.Cells.Item(79, 4).Value2 = 'home'
$wsUnit.Cells.Item(79, 5).Value2 = 1

$wsUnit.Cells.Item(80, 1).Value2 = 'Meltingar- og nýrnadeild (Hb-12E)'
$wsUnit.Cells.Item(80, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(80, 3).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(80, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(80, 5).Value2 = 2

$wsUnit.Cells.Item(81, 1).Value2 = 'Vökudeild (Hb-23D)'
$wsUnit.Cells.Item(81, 2).Value2 = 'Legudeild'
$wsUnit.Cells.Item(81, 3).Value2 = 'inpatient_ward_pediadric'
$wsUnit.Cells.Item(81, 4).Value2 = 'inpatient_ward'
$wsUnit.Cells.Item(81, 5).Value2 = 2

# --- Column C width adjustment on lsh_unit_categories (best achievable approximation) ---
$wsUnit.Columns.Item(3).ColumnWidth = 27.8

# --- View/selection state ---
$wsIsolation = $wb.Worksheets.Item("lsh_isolation_categories")
[void]$wsIsolation.Activate()
[void]$wsIsolation.Range("B21").Select()

[void]$wsUnit.Activate()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
[void]$wsUnit.Range("C85").Select()

Write-Output "Edit complete"
